$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 691-692; existing rows 691-720 shift down to 693-722.
$ws.Range("A691:A692").EntireRow.Insert()

# Populate the two newly-inserted rows with the new weekly price records.
$ws.Cells.Item(691, 1).Value = 9
$ws.Cells.Item(691, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(691, 3).Value = "Metropolitana"
$ws.Cells.Item(691, 4).Value = 44509
$ws.Cells.Item(691, 5).Value = 13
$ws.Cells.Item(691, 6).Value = 100114001
$ws.Cells.Item(691, 7).Value = "Papa"
$ws.Cells.Item(691, 8).Value = "Asterix"
$ws.Cells.Item(691, 9).Value = "1a (nueva lavada)"
$ws.Cells.Item(691, 10).Value = 340
$ws.Cells.Item(691, 11).Value = 11000
$ws.Cells.Item(691, 12).Value = 12000
$ws.Cells.Item(691, 13).Value = 11500
$ws.Cells.Item(691, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(691, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(691, 16).Value = 460
$ws.Cells.Item(691, 17).Value = 25
$ws.Cells.Item(691, 18).Value = "Hortaliza"

$ws.Cells.Item(692, 1).Value = 9
$ws.Cells.Item(692, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(692, 3).Value = "Metropolitana"
$ws.Cells.Item(692, 4).Value = 44509
$ws.Cells.Item(692, 5).Value = 13
$ws.Cells.Item(692, 6).Value = 100114001
$ws.Cells.Item(692, 7).Value = "Papa"
$ws.Cells.Item(692, 8).Value = "Asterix"
$ws.Cells.Item(692, 9).Value = "1a nueva(o)"
$ws.Cells.Item(692, 10).Value = 196
$ws.Cells.Item(692, 11).Value = 9000
$ws.Cells.Item(692, 12).Value = 10000
$ws.Cells.Item(692, 13).Value = 9500
$ws.Cells.Item(692, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(692, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(692, 16).Value = 380
$ws.Cells.Item(692, 17).Value = 25
$ws.Cells.Item(692, 18).Value = "Hortaliza"
